$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update a few existing "Underscore" group test cases (rows 12-14) ---
$ws.Range("C12").Value = 'Select BOOLEAN pattern a(_, _"2+20*a"_)'
$ws.Range("C13").Value = 'Select BOOLEAN pattern a(_, _"(f+3)"_)'
$ws.Range("C14").Value = 'Select BOOLEAN pattern a(_, _"(a+b)"_)'

# --- Insert 3 new rows after row 16 for new "Not WellFormed" BOOLEAN test cases ---
$ws.Rows("17:19").Insert()

$ws.Range("A17").Value = 16
$ws.Range("B17").Value = "assign a;"
$ws.Range("C17").Value = 'Select BOOLEAN pattern a(_, _"()"_)'
$ws.Range("D17").Value = "false"
$ws.Range("E17").Value = "Not WellFormed Expr BracketDontMatch - Test Valdiation of Expression"

$ws.Range("A18").Value = 17
$ws.Range("B18").Value = "assign a;"
$ws.Range("C18").Value = 'Select BOOLEAN pattern a(_, _"+4"_)'
$ws.Range("D18").Value = "false"
$ws.Range("E18").Value = "Not WellFormed Expr BracketDontMatch - Test Valdiation of Expression"

$ws.Range("A19").Value = 18
$ws.Range("B19").Value = "assign a;"
$ws.Range("C19").Value = 'Select BOOLEAN pattern a(_, _"1-"_)'
$ws.Range("D19").Value = "false"
$ws.Range("E19").Value = "Not WellFormed Expr BracketDontMatch - Test Valdiation of Expression"

# --- Re-number the Index column for all rows pushed down by the insert ---
$ws.Range("A20").Value = 19
$ws.Range("A21").Value = 20
$ws.Range("A22").Value = 21
$ws.Range("A23").Value = 22
$ws.Range("A24").Value = 23
$ws.Range("A25").Value = 24
$ws.Range("A26").Value = 25
$ws.Range("A27").Value = 26
$ws.Range("A28").Value = 27
$ws.Range("A29").Value = 28
$ws.Range("A30").Value = 29

# --- Update the Ident PartialMatch Expr test case (now at row 23) ---
$ws.Range("C23").Value = 'Select BOOLEAN pattern a("f", _"c * b + 2- 5"_)'

# --- Drop "such that " from the remaining query texts (now rows 25-30) ---
$ws.Range("C25").Value = "Select v pattern a(v, _)"
$ws.Range("C26").Value = 'Select v pattern a(v, _"b"_)'
$ws.Range("C27").Value = 'Select v pattern a(v, "2")'
$ws.Range("C28").Value = "Select a pattern a(v, _)"
$ws.Range("C29").Value = 'Select a pattern a(v, _"b"_)'
$ws.Range("C30").Value = 'Select a pattern a(v, "2")'

# --- Update the frozen-pane / selection view state to match the new extent ---
$ws.Range("C30").Select()
